$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": MyForecast for week W7 (row 8) revised up ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D8").Value = 25

# --- Sheet "Summary": update derived totals / max-forecast week              ---
# These cells are stored as literal text (e.g. "336", "2025-03-16") rather
# than numbers/dates, so a plain Range.Value assignment would auto-coerce
# numeric-looking / date-looking text to a Number/Date (exactly like typing
# such a string into a General-formatted cell in Excel), which would not
# reproduce the original text type. Instead, a scratch cell is loaded with
# a `="<text>"` formula (always a text result), copied, and paste-special'd
# as values-only into the destination - this writes a genuine text value
# without touching the destination cell's number format/style.
$wsSummary = $wb.Worksheets.Item("Summary")
$scratch = $wsSummary.Range("Z1")

function Set-LiteralText($range, [string]$text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

Set-LiteralText $wsSummary.Range("B9") "345"
Set-LiteralText $wsSummary.Range("B10") "163"
Set-LiteralText $wsSummary.Range("B13") "2025-03-09"

$scratch.ClearContents()
$excel.CutCopyMode = 0
